$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B-E keep their original text formatting while we overwrite values
# (some new values look numeric, e.g. "606.51", and must stay as text to match the source data)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.556.13'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.520.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.51'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.84'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.517.16'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.196'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.23'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +8.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.588'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.43'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.67%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.091.64'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '612.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.529.50'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.625.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.93%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.883'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -9.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.62'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.58'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.07'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.86%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.30%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.98'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '641.93'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +12.63%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0999'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.57'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0478'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +6.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.86'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0750'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.367.19'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.85%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '32.25'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.02%  '
